# Auto-generated edit script: apply literal value updates to match target diff.
# All target cells are plain numeric literals (no formulas in the workbook).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2301
$ws.Range("I19").Value = 2301
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2301
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -2126
$ws.Range("N19").ClearContents()
$ws.Range("H43").Value = 651
$ws.Range("J43").Value = 651
$ws.Range("L43").Value = 651
$ws.Range("N43").Value = -789
$ws.Range("H58").Value = 69803.55499999999
$ws.Range("I58").Value = 166743
$ws.Range("J58").Value = 21333.834
$ws.Range("K58").Value = 500229
$ws.Range("L58").Value = 64001.50199999999
$ws.Range("M58").Value = -500079
$ws.Range("N58").Value = -64301.50199999999
$ws.Range("H96").Value = 998
$ws.Range("J96").Value = 998
$ws.Range("L96").Value = 2994
$ws.Range("N96").Value = -5740
$ws.Range("H116").Value = 27785944
$ws.Range("J116").Value = 11876.25
$ws.Range("L116").Value = 11876.25
$ws.Range("N116").Value = -18760.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2968.7144
$ws.Range("I2").Value = 1258.6666
$ws.Range("J2").Value = 5248.778
$ws.Range("K2").Value = 1258.6666
$ws.Range("L2").Value = 5248.778
$ws.Range("M2").Value = -1145.6666
$ws.Range("N2").Value = -5474.778
$ws.Range("H32").Value = 5656.125
$ws.Range("I32").Value = 5065.7915
$ws.Range("J32").Value = 9198.125
$ws.Range("K32").Value = 5065.7915
$ws.Range("L32").Value = 9198.125
$ws.Range("M32").Value = -4778.7915
$ws.Range("N32").Value = -9772.125
$ws.Range("H57").Value = 4552.4287
$ws.Range("I57").Value = 4552.4287
$ws.Range("K57").Value = 4552.4287
$ws.Range("M57").Value = -4068.4287
$ws.Range("H116").Value = 2968.7144
$ws.Range("I116").Value = 1258.6666
$ws.Range("J116").Value = 5248.778
$ws.Range("K116").Value = 1258.6666
$ws.Range("L116").Value = 5248.778
$ws.Range("M116").Value = 1035.3334
$ws.Range("N116").Value = -9836.778
$ws.Range("H132").Value = 4569.8306
$ws.Range("I132").Value = 3552.0212
$ws.Range("J132").Value = 7227.4443
$ws.Range("K132").Value = 10656.0636
$ws.Range("L132").Value = 21682.3329
$ws.Range("M132").Value = -8126.063600000001
$ws.Range("N132").Value = -26742.3329

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2968.7144
$ws.Range("I3").Value = 1258.6666
$ws.Range("J3").Value = 5248.778
$ws.Range("K3").Value = 1258.6666
$ws.Range("L3").Value = 5248.778
$ws.Range("M3").Value = -1144.6666
$ws.Range("N3").Value = -5476.778
$ws.Range("H22").Value = 259.4
$ws.Range("I22").Value = 232.33333
$ws.Range("K22").Value = 232.33333
$ws.Range("M22").Value = -59.33332999999999
$ws.Range("H134").Value = 4573.1787
$ws.Range("I134").Value = 1702.8889
$ws.Range("K134").Value = 5108.6667
$ws.Range("M134").Value = -2573.6667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 22737038
$ws.Range("I58").Value = 83335830
$ws.Range("J58").Value = 12489.9375
$ws.Range("K58").Value = 83335830
$ws.Range("L58").Value = 12489.9375
$ws.Range("M58").Value = -83335627
$ws.Range("N58").Value = -12895.9375
$ws.Range("H68").Value = 85000
$ws.Range("J68").Value = 85000
$ws.Range("L68").Value = 85000
$ws.Range("N68").Value = -86498
$ws.Range("H71").Value = 85000
$ws.Range("J71").Value = 85000
$ws.Range("L71").Value = 255000
$ws.Range("N71").Value = -262488
$ws.Range("H132").Value = 4619.0684
$ws.Range("I132").Value = 1478.3334
$ws.Range("K132").Value = 4435.0002
$ws.Range("M132").Value = -1905.0002
$ws.Range("H136").Value = 22737038
$ws.Range("I136").Value = 83335830
$ws.Range("J136").Value = 12489.9375
$ws.Range("K136").Value = 250007490
$ws.Range("L136").Value = 37469.8125
$ws.Range("M136").Value = -250004940
$ws.Range("N136").Value = -42569.8125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 453.66666
$ws.Range("I86").Value = 454.0909
$ws.Range("J86").Value = 449
$ws.Range("K86").Value = 1362.2727
$ws.Range("L86").Value = 1347
$ws.Range("M86").Value = -176.2727
$ws.Range("N86").Value = -3719
$ws.Range("H89").Value = 453.66666
$ws.Range("I89").Value = 454.0909
$ws.Range("J89").Value = 449
$ws.Range("K89").Value = 4086.8181
$ws.Range("L89").Value = 4041
$ws.Range("M89").Value = 1841.1819
$ws.Range("N89").Value = -15897
$ws.Range("H92").Value = 7693730.5
$ws.Range("I92").Value = 32
$ws.Range("J92").Value = 8548586
$ws.Range("K92").Value = 96
$ws.Range("L92").Value = 25645758
$ws.Range("M92").Value = 1152
$ws.Range("N92").Value = -25648254
$ws.Range("H122").Value = 1888471.6
$ws.Range("I122").Value = 4715684.5
$ws.Range("K122").Value = 42441160.5
$ws.Range("M122").Value = -42438710.5
$ws.Range("H132").Value = 8877.925999999999
$ws.Range("I132").Value = 4934.8184
$ws.Range("J132").Value = 11588.8125
$ws.Range("K132").Value = 44413.3656
$ws.Range("L132").Value = 104299.3125
$ws.Range("M132").Value = -41883.3656
$ws.Range("N132").Value = -109359.3125
$ws.Range("H138").Value = 4733.55
$ws.Range("I138").Value = 3727.75
$ws.Range("K138").Value = 11183.25
$ws.Range("M138").Value = -6043.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3149.75
$ws.Range("I107").Value = 3533
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 3533
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -1613
$ws.Range("N107").Value = -5840
$ws.Range("H121").Value = 56546
$ws.Range("J121").Value = 56546
$ws.Range("L121").Value = 56546
$ws.Range("N121").Value = -60040

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 2422
$ws.Range("I58").Value = 2422
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2422
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2162
$ws.Range("N58").ClearContents()
$ws.Range("H59").Value = 56447.332
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 56447.332
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 56447.332
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -57755.332
$ws.Range("H105").Value = 60597
$ws.Range("J105").Value = 60597
$ws.Range("L105").Value = 60597
$ws.Range("N105").Value = -67585
$ws.Range("H106").Value = 39418.668
$ws.Range("J106").Value = 39418.668
$ws.Range("L106").Value = 39418.668
$ws.Range("N106").Value = -41942.668
$ws.Range("H110").Value = 46167
$ws.Range("J110").Value = 46167
$ws.Range("L110").Value = 46167
$ws.Range("N110").Value = -54347
$ws.Range("H120").Value = 48479
$ws.Range("J120").Value = 48479
$ws.Range("L120").Value = 48479
$ws.Range("N120").Value = -58155
$ws.Range("H122").Value = 5016.077
$ws.Range("I122").Value = 2894.8
$ws.Range("J122").Value = 7908.727
$ws.Range("K122").Value = 8684.400000000001
$ws.Range("L122").Value = 23726.181
$ws.Range("M122").Value = -6234.400000000001
$ws.Range("N122").Value = -28626.181
$ws.Range("H125").Value = 51598
$ws.Range("J125").Value = 51598
$ws.Range("L125").Value = 51598
$ws.Range("N125").Value = -61438
$ws.Range("H127").Value = 59519
$ws.Range("J127").Value = 59519
$ws.Range("L127").Value = 59519
$ws.Range("N127").Value = -69439
$ws.Range("H132").Value = 10422844
$ws.Range("I132").Value = 22729432
$ws.Range("K132").Value = 68188296
$ws.Range("M132").Value = -68185766
$ws.Range("H136").Value = 11912.652

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1602.5769
$ws.Range("I113").Value = 722.93335
$ws.Range("J113").Value = 2802.0908
$ws.Range("K113").Value = 2168.80005
$ws.Range("L113").Value = 8406.2724
$ws.Range("M113").Value = 1.199950000000172
$ws.Range("N113").Value = -12746.2724
$ws.Range("H126").Value = 4614.4
$ws.Range("I126").Value = 3236.25
$ws.Range("J126").Value = 5533.1665
$ws.Range("K126").Value = 9708.75
$ws.Range("L126").Value = 16599.4995
$ws.Range("M126").Value = -7238.75
$ws.Range("N126").Value = -21539.4995
$ws.Range("H135").Value = 71811
$ws.Range("J135").Value = 71811
$ws.Range("L135").Value = 71811
$ws.Range("N135").Value = -81951

